$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.181.28'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.801.40'
$ws.Range("E3").Value = '  +4.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '618.71'
$ws.Range("E5").Value = '  +4.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.26'
$ws.Range("E6").Value = '  -3.67%  '
$ws.Range("D7").Value = '3.801.82'
$ws.Range("E7").Value = '  +4.26%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("E10").Value = '  +4.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.35'
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.495'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.12'
$ws.Range("E13").Value = '  +4.67%  '
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").Value = '4.435.10'
$ws.Range("E15").Value = '  +4.68%  '
$ws.Range("D16").Value = '3.800.26'
$ws.Range("E16").Value = '  +4.80%  '
$ws.Range("D17").Value = '70.191.03'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("E19").Value = '  +1.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '515.22'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.68'
$ws.Range("E21").Value = '  -3.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.61'
$ws.Range("E22").Value = '  +4.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.730'
$ws.Range("E23").Value = '  -2.82%  '
$ws.Range("E24").Value = '  +5.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.21'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.30'
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.22'
$ws.Range("E27").Value = '  +3.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000139'
$ws.Range("E28").Value = '  +27.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.49'
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.84'
$ws.Range("E31").Value = '  +3.15%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.82'
$ws.Range("E32").Value = '  -5.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.84'
$ws.Range("E33").Value = '  -0.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.115'
$ws.Range("E34").Value = '  -2.01%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.22'
$ws.Range("E36").Value = '  +1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").Value = '  +3.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.340'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.17'
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  +3.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.41'
$ws.Range("E41").Value = '  +1.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '44.37'
$ws.Range("E42").Value = '  -5.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.80'
$ws.Range("E43").Value = '  -1.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '422.78'
$ws.Range("E44").Value = '  +5.05%  '
$ws.Range("D45").Value = '3.064.76'
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.77'
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0367'
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.57'
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '136.36'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.49'
$ws.Range("E51").Value = '  +0.88%  '
